$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 60, pushing the existing rows 60-68 down to 61-69
$ws.Rows.Item(60).Insert()

# Populate the newly inserted row 60 with the new weekly price observation
$ws.Cells.Item(60, 1).Value = 3
$ws.Cells.Item(60, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(60, 3).Value = "Coquimbo"
$ws.Cells.Item(60, 4).Value = 44476
$ws.Cells.Item(60, 5).Value = 5
$ws.Cells.Item(60, 6).Value = 100112026
$ws.Cells.Item(60, 7).Value = "Haba"
$ws.Cells.Item(60, 8).Value = "Sin especificar"
$ws.Cells.Item(60, 9).Value = "Primera"
$ws.Cells.Item(60, 10).Value = 73
$ws.Cells.Item(60, 11).Value = 8500
$ws.Cells.Item(60, 12).Value = 9000
$ws.Cells.Item(60, 13).Value = 8740
$ws.Cells.Item(60, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(60, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(60, 16).Value = 350
$ws.Cells.Item(60, 17).Value = 25
$ws.Cells.Item(60, 18).Value = "Hortaliza"
